$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n (column J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary labels + stats
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold/size-12/vertical-centered style for the new stat values (build once, then copy the format)
$cell = $ws.Range("B14")
$cell.Font.Bold = $true
$cell.Font.Size = 12
$cell.VerticalAlignment = -4108

$cell.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the new label rows
$ws.Range("A14:B17").RowHeight = 15.6

# Selection left behind by the edit
$ws.Range("A14:B17").Select()

# Print setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
